$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.915.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.324.22"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -5.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.60"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "84.25"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -8.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.530"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.47%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.482"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0811"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "29.91"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -8.76%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.676.74"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -5.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.40"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.61"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.321.67"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.751"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "39.866.41"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0896"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.06"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.46"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.57"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.88"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.54"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -7.40%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.79"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -7.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.10"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -6.64%  "
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.23"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.53"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "150.33"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.59%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.08"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.31%  "
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0717"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -6.10%  "
$ws.Range("E36").Value = "  -2.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0992"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.74"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.51"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -8.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.70"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -6.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.79"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.02%  "
$ws.Range("E42").Value = "  -3.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.944.76"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0264"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.43"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.36"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.67"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -9.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.539.83"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.37"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.13"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.23"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.32%  "
